$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 281-282, pushing the existing rows 281..328 down to 283..330.
$ws.Rows("281:282").Insert()

# New row 281: Terminal La Palmera de La Serena, Pera, Packham's Triumph, Primera
$ws.Range("A281").Value = 8
$ws.Range("B281").Value = "Terminal La Palmera de La Serena"
$ws.Range("C281").Value = "Coquimbo"
$ws.Range("D281").Value = 44474
$ws.Range("E281").Value = 4
$ws.Range("F281").Value = "Fruta"
$ws.Range("G281").Value = 100104
$ws.Range("H281").Value = "Frutos de pepita"
$ws.Range("I281").Value = 100104005
$ws.Range("J281").Value = "Pera"
$ws.Range("K281").Value = "Packham's Triumph"
$ws.Range("L281").Value = "Primera"
$ws.Range("M281").Value = 24
$ws.Range("N281").Value = 240000
$ws.Range("O281").Value = 250000
$ws.Range("P281").Value = 245000
$ws.Range("Q281").Value = "`$/bins (450 kilos)"
$ws.Range("R281").Value = "Región de O'Higgins"
$ws.Range("S281").Value = 544
$ws.Range("T281").Value = 450

# New row 282: Terminal La Palmera de La Serena, Pera, Packham's Triumph, Segunda
$ws.Range("A282").Value = 8
$ws.Range("B282").Value = "Terminal La Palmera de La Serena"
$ws.Range("C282").Value = "Coquimbo"
$ws.Range("D282").Value = 44474
$ws.Range("E282").Value = 4
$ws.Range("F282").Value = "Fruta"
$ws.Range("G282").Value = 100104
$ws.Range("H282").Value = "Frutos de pepita"
$ws.Range("I282").Value = 100104005
$ws.Range("J282").Value = "Pera"
$ws.Range("K282").Value = "Packham's Triumph"
$ws.Range("L282").Value = "Segunda"
$ws.Range("M282").Value = 18
$ws.Range("N282").Value = 210000
$ws.Range("O282").Value = 220000
$ws.Range("P282").Value = 215000
$ws.Range("Q282").Value = "`$/bins (450 kilos)"
$ws.Range("R282").Value = "Región de O'Higgins"
$ws.Range("S282").Value = 478
$ws.Range("T282").Value = 450
